$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 71.333336
$ws.Range("I4").Value = 71.333336
$ws.Range("K4").Value = 71.333336
$ws.Range("M4").Value = 42.666664

$ws.Range("H41").Value = 342.72726
$ws.Range("I41").Value = 397.1111
$ws.Range("J41").Value = 98
$ws.Range("K41").Value = 397.1111
$ws.Range("L41").Value = 98
$ws.Range("M41").Value = 42.88889999999998
$ws.Range("N41").Value = -978

$ws.Range("H64").Value = 7600
$ws.Range("J64").Value = 7600
$ws.Range("L64").Value = 7600
$ws.Range("N64").Value = -8096

$ws.Range("H67").Value = 7600
$ws.Range("J67").Value = 7600
$ws.Range("L67").Value = 7600
$ws.Range("N67").Value = -9316

$ws.Range("H69").Value = 7614.423
$ws.Range("J69").Value = 7614.423
$ws.Range("L69").Value = 22843.269
$ws.Range("N69").Value = -24591.269

$ws.Range("H72").Value = 7614.423
$ws.Range("J72").Value = 7614.423
$ws.Range("L72").Value = 68529.807
$ws.Range("N72").Value = -77265.807

$ws.Range("H88").Value = 916.5714
$ws.Range("I88").Value = 1722
$ws.Range("J88").Value = 594.4
$ws.Range("K88").Value = 1722
$ws.Range("L88").Value = 594.4
$ws.Range("M88").Value = -1316
$ws.Range("N88").Value = -1406.4

$ws.Range("H91").Value = 916.5714
$ws.Range("I91").Value = 1722
$ws.Range("J91").Value = 594.4
$ws.Range("K91").Value = 1722
$ws.Range("L91").Value = 594.4
$ws.Range("M91").Value = -318
$ws.Range("N91").Value = -3402.4

$ws.Range("H100").Value = 1650.5
$ws.Range("I100").Value = 1680.6
$ws.Range("K100").Value = 1680.6
$ws.Range("M100").Value = -1139.6

$ws.Range("H116").Value = 4714.3335
$ws.Range("I116").Value = 4719.778
$ws.Range("J116").Value = 4698
$ws.Range("K116").Value = 4719.778
$ws.Range("L116").Value = 4698
$ws.Range("M116").Value = -1277.778
$ws.Range("N116").Value = -11582

$ws.Range("H125").Value = 950
$ws.Range("I125").Value = 950
$ws.Range("K125").Value = 8550
$ws.Range("M125").Value = -6090

$ws.Range("H132").Value = 3374.2432
$ws.Range("I132").Value = 1188.6129
$ws.Range("K132").Value = 3565.8387
$ws.Range("M132").Value = -1035.8387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = $null

$ws.Range("H45").Value = 3460
$ws.Range("I45").Value = 1971
$ws.Range("J45").Value = 4390.625
$ws.Range("K45").Value = 1971
$ws.Range("L45").Value = 4390.625
$ws.Range("M45").Value = -1594
$ws.Range("N45").Value = -5144.625

$ws.Range("H74").Value = 3409.5454
$ws.Range("I74").Value = 2653.4443
$ws.Range("K74").Value = 2653.4443
$ws.Range("M74").Value = -1779.4443

$ws.Range("H77").Value = 3409.5454
$ws.Range("I77").Value = 2653.4443
$ws.Range("K77").Value = 13267.2215
$ws.Range("M77").Value = -8899.2215

$ws.Range("H96").Value = 4020088.8
$ws.Range("J96").Value = 4020088.8
$ws.Range("L96").Value = 4020088.8
$ws.Range("N96").Value = -4025580.8

$ws.Range("H97").Value = 1116.3572
$ws.Range("I97").Value = 1212.2727
$ws.Range("J97").Value = 764.6667
$ws.Range("K97").Value = 1212.2727
$ws.Range("L97").Value = 764.6667
$ws.Range("M97").Value = -716.2727
$ws.Range("N97").Value = -1756.6667

$ws.Range("H122").Value = 2636.6667
$ws.Range("I122").Value = 2636.6667
$ws.Range("K122").Value = 7910.000100000001
$ws.Range("M122").Value = -5460.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4010
$ws.Range("I99").Value = 4010
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4010
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2512
$ws.Range("N99").Value = $null

$ws.Range("H105").Value = 1190.7
$ws.Range("I105").Value = 1190.7
$ws.Range("K105").Value = 1190.7
$ws.Range("M105").Value = 556.3

$ws.Range("H107").Value = 9200.286
$ws.Range("I107").Value = 9235.333000000001
$ws.Range("K107").Value = 9235.333000000001
$ws.Range("M107").Value = -7315.333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4397.8623
$ws.Range("I31").Value = 1269.2667
$ws.Range("J31").Value = 7749.9287
$ws.Range("K31").Value = 1269.2667
$ws.Range("L31").Value = 7749.9287
$ws.Range("M31").Value = -974.2666999999999
$ws.Range("N31").Value = -8339.9287

$ws.Range("H34").Value = 4397.8623
$ws.Range("I34").Value = 1269.2667
$ws.Range("J34").Value = 7749.9287
$ws.Range("K34").Value = 1269.2667
$ws.Range("L34").Value = 7749.9287
$ws.Range("M34").Value = -1067.2667
$ws.Range("N34").Value = -8153.9287

$ws.Range("H58").Value = 2890.2
$ws.Range("I58").Value = 1451.9
$ws.Range("K58").Value = 1451.9
$ws.Range("M58").Value = -1248.9

$ws.Range("H135").Value = 82854.5
$ws.Range("I135").Value = 80709
$ws.Range("J135").Value = 85000
$ws.Range("K135").Value = 80709
$ws.Range("L135").Value = 85000
$ws.Range("M135").Value = -75639
$ws.Range("N135").Value = -95140

$ws.Range("H136").Value = 2890.2
$ws.Range("I136").Value = 1451.9
$ws.Range("K136").Value = 4355.700000000001
$ws.Range("M136").Value = -1805.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3397.4
$ws.Range("I132").Value = 1494
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 13446
$ws.Range("L132").Value = 41997.0015
$ws.Range("M132").Value = -10916
$ws.Range("N132").Value = -47057.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 337.5
$ws.Range("I13").Value = 50
$ws.Range("J13").Value = 625
$ws.Range("K13").Value = 50
$ws.Range("L13").Value = 625
$ws.Range("M13").Value = 89
$ws.Range("N13").Value = -903

$ws.Range("H20").Value = 41269.43
$ws.Range("J20").Value = 41269.43
$ws.Range("L20").Value = 41269.43
$ws.Range("N20").Value = -41759.43

$ws.Range("H24").Value = 38094.855
$ws.Range("J24").Value = 38094.855
$ws.Range("L24").Value = 38094.855
$ws.Range("N24").Value = -38440.855

$ws.Range("H80").Value = 1132.5
$ws.Range("J80").Value = 1119.2
$ws.Range("L80").Value = 1119.2
$ws.Range("N80").Value = -3115.2

$ws.Range("H83").Value = 1132.5
$ws.Range("J83").Value = 1119.2
$ws.Range("L83").Value = 5596
$ws.Range("N83").Value = -15580

$ws.Range("H102").Value = 3438.9285
$ws.Range("I102").Value = 3170.4167
$ws.Range("J102").Value = 5050
$ws.Range("K102").Value = 3170.4167
$ws.Range("L102").Value = 5050
$ws.Range("M102").Value = -1548.4167
$ws.Range("N102").Value = -8294

$ws.Range("H113").Value = 9580
$ws.Range("I113").Value = 8000
$ws.Range("K113").Value = 8000
$ws.Range("M113").Value = -5830

$ws.Range("H132").Value = 168766.83
$ws.Range("I132").Value = 168766.83
$ws.Range("K132").Value = 506300.49
$ws.Range("M132").Value = -503770.49

$ws.Range("H133").Value = 109999.664
$ws.Range("J133").Value = 109999.664
$ws.Range("L133").Value = 109999.664
$ws.Range("N133").Value = -120119.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1163.7
$ws.Range("I93").Value = 1148.1428
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 1148.1428
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = 99.85719999999992
$ws.Range("N93").Value = -3696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336

$ws.Range("H31").Value = 12062.5
$ws.Range("I31").Value = 8583.333000000001
$ws.Range("J31").Value = 22500
$ws.Range("K31").Value = 8583.333000000001
$ws.Range("L31").Value = 22500
$ws.Range("M31").Value = -8235.333000000001
$ws.Range("N31").Value = -23196

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null

$ws.Range("H62").Value = 9541.166999999999
$ws.Range("I62").Value = 7249.25
$ws.Range("K62").Value = 7249.25
$ws.Range("M62").Value = -6625.25

$ws.Range("H65").Value = 9541.166999999999
$ws.Range("I65").Value = 7249.25
$ws.Range("K65").Value = 36246.25
$ws.Range("M65").Value = -33126.25

$ws.Range("H81").Value = 537.5
$ws.Range("I81").Value = 537.5
$ws.Range("K81").Value = 1075
$ws.Range("M81").Value = -14

$ws.Range("H84").Value = 537.5
$ws.Range("I84").Value = 537.5
$ws.Range("K84").Value = 5375
$ws.Range("M84").Value = -71

$ws.Range("H96").Value = 1165.3334
$ws.Range("I96").Value = 1096
$ws.Range("K96").Value = 1096
$ws.Range("M96").Value = 277

$ws.Range("H132").Value = 2477.2856
$ws.Range("I132").Value = 1968.6
$ws.Range("K132").Value = 5905.799999999999
$ws.Range("M132").Value = -3375.799999999999

$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120

$ws.Range("H139").Value = 74266.664
$ws.Range("J139").Value = 74266.664
$ws.Range("L139").Value = 74266.664
$ws.Range("N139").Value = -84546.664

$ws.Range("H141").Value = 416473.34
$ws.Range("I141").Value = 500000
$ws.Range("J141").Value = 374710
$ws.Range("K141").Value = 500000
$ws.Range("L141").Value = 374710
$ws.Range("M141").Value = -494820
$ws.Range("N141").Value = -385070

Write-Output "Applied all changes"